# Auto-generated edit script: updates the cryptos price/volume table
# to match the refreshed scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.875.07'
$ws.Range('E2').Value = '  -4.19%  '
$ws.Range('D3').Value = '2.993.78'
$ws.Range('E3').Value = '  -4.07%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '544.26'
$ws.Range('E5').Value = '  -4.80%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '152.33'
$ws.Range('E6').Value = '  -6.51%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.571'
$ws.Range('E8').Value = '  -0.22%  '
$ws.Range('D9').Value = '3.006.59'
$ws.Range('E9').Value = '  -3.76%  '
$ws.Range('E10').Value = '  -2.97%  '
$ws.Range('E11').Value = '  -6.94%  '
$ws.Range('E12').Value = '  -2.81%  '
$ws.Range('D13').Value = '3.517.99'
$ws.Range('E13').Value = '  -3.94%  '
$ws.Range('E14').Value = '  -0.69%  '
$ws.Range('D15').Value = '61.872.09'
$ws.Range('E15').Value = '  -4.22%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '24.08'
$ws.Range('E16').Value = '  -2.86%  '
$ws.Range('D17').Value = '2.998.74'
$ws.Range('E17').Value = '  -4.24%  '
$ws.Range('E18').Value = '  -4.69%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.19'
$ws.Range('E19').Value = '  -0.60%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '12.08'
$ws.Range('E20').Value = '  -2.49%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '379.72'
$ws.Range('E21').Value = '  -7.03%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '6.74'
$ws.Range('E22').Value = '  -3.69%  '
$ws.Range('E23').Value = '  +0.00%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '66.16'
$ws.Range('E24').Value = '  -2.72%  '
$ws.Range('D25').Value = '3.116.52'
$ws.Range('E25').Value = '  -4.76%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '0.471'
$ws.Range('E26').Value = '  -1.77%  '
$ws.Range('E27').Value = '  -2.45%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.997'
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('D29').Value = '0.0₃0937'
$ws.Range('E29').Value = '  -8.53%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.28'
$ws.Range('E30').Value = '  -9.30%  '
$ws.Range('E31').Value = '  +0.01%  '
$ws.Range('E32').Value = '  -4.10%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '20.50'
$ws.Range('E33').Value = '  -3.40%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '160.44'
$ws.Range('E34').Value = '  -2.33%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '4.61'
$ws.Range('E35').Value = '  -6.34%  '
$ws.Range('B36').Value = 'Aptos'
$ws.Range('C36').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '5.95'
$ws.Range('E36').Value = '  -4.22%  '
$ws.Range('E37').Value = '  -5.29%  '
$ws.Range('E38').Value = '  -5.19%  '
$ws.Range('E39').Value = '  -6.65%  '
$ws.Range('B40').Value = 'Filecoin'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.91'
$ws.Range('E40').Value = '  -4.71%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '37.52'
$ws.Range('E41').Value = '  -2.03%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.419.16'
$ws.Range('E42').Value = '  -6.87%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '22.14'
$ws.Range('E43').Value = '  -6.27%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.674'
$ws.Range('E44').Value = '  -2.00%  '
$ws.Range('E45').Value = '  -3.78%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '5.26'
$ws.Range('E46').Value = '  +0.60%  '
$ws.Range('E47').Value = '  +0.10%  '
$ws.Range('E48').Value = '  -3.56%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0955'
$ws.Range('E49').Value = '  -1.73%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '269.14'
$ws.Range('E50').Value = '  -5.84%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '19.69'
$ws.Range('E51').Value = '  -6.46%  '
